# Auto commit at 2025-12-07  8:21:21.00
# Updates the monthly "Metrics" figures (charging kwh / revenue / orders, etc.)
# to the latest pull. Values on the "today" sheet reference Metrics!B2:B13
# via formulas and recalculate automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 81509.62
$metrics.Range("B3").Value  = 70078.33
$metrics.Range("B4").Value  = 25111.809999999998
$metrics.Range("B5").Value  = 3321
$metrics.Range("B6").Value  = 5284216.7300000004
$metrics.Range("B7").Value  = 4470431.290000001
$metrics.Range("B8").Value  = 1557068.6900000004
$metrics.Range("B9").Value  = 206028
$metrics.Range("B10").Value = 33749597.719999999
$metrics.Range("B11").Value = 31745706.450000003
$metrics.Range("B12").Value = 11838790.729999997
$metrics.Range("B13").Value = 1303658

# Restore the selection cursor positions recorded by the author at save time.
[void]$metrics.Range("E28").Select()

$today = $wb.Worksheets.Item("today")
[void]$today.Range("E10").Select()
